$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.248.83'
$ws.Range('E2').Value = '  -4.06%  '
$ws.Range('D3').Value = '1.657.59'
$ws.Range('E3').Value = '  -3.47%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '216.73'
$ws.Range('E5').Value = '  -3.68%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5149'
$ws.Range('E6').Value = '  -3.02%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.008'
$ws.Range('E7').Value = '  +0.30%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2604'
$ws.Range('E8').Value = '  -2.11%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06468'
$ws.Range('E9').Value = '  -3.74%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.98'
$ws.Range('E10').Value = '  -4.71%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07843'
$ws.Range('E11').Value = '  +1.79%  '
$ws.Range('D12').Value = '1.664.49'
$ws.Range('E12').Value = '  -3.06%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.315'
$ws.Range('E13').Value = '  -3.98%  '
$ws.Range('D14').Value = '1.886.28'
$ws.Range('E14').Value = '  -3.46%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.5547'
$ws.Range('E15').Value = '  -4.76%  '
$ws.Range('D16').Value = '0.0₅8053'
$ws.Range('E16').Value = '  -2.13%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '64.37'
$ws.Range('E17').Value = '  -5.35%  '
$ws.Range('D18').Value = '26.257.29'
$ws.Range('E18').Value = '  -4.14%  '
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '210.89'
$ws.Range('E20').Value = '  -5.54%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.423'
$ws.Range('E21').Value = '  -5.30%  '
$ws.Range('E22').Value = '  -3.45%  '
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.008'
$ws.Range('E24').Value = '  +0.27%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '144.81'
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('E26').Value = '  +3.56%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.1176'
$ws.Range('E27').Value = '  -2.78%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.009'
$ws.Range('E28').Value = '  -3.42%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.87'
$ws.Range('E29').Value = '  -2.59%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05118'
$ws.Range('E30').Value = '  -5.43%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.248'
$ws.Range('E31').Value = '  -3.78%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.368'
$ws.Range('E32').Value = '  -3.48%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.237'
$ws.Range('E33').Value = '  -5.39%  '
$ws.Range('E34').Value = '  -4.25%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.740'
$ws.Range('E35').Value = '  -4.31%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.9279'
$ws.Range('E36').Value = '  -2.94%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.358'
$ws.Range('E37').Value = '  -1.47%  '
$ws.Range('D38').Value = '1.175.06'
$ws.Range('E38').Value = '  +2.41%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.5744'
$ws.Range('E39').Value = '  -2.89%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01592'
$ws.Range('E40').Value = '  -3.94%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.564'
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.008'
$ws.Range('E42').Value = '  +0.27%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.722'
$ws.Range('E43').Value = '  -1.93%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.8261'
$ws.Range('E44').Value = '  -1.88%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '100.52'
$ws.Range('E45').Value = '  -0.60%  '
$ws.Range('D46').Value = '1.797.52'
$ws.Range('E46').Value = '  -3.39%  '
$ws.Range('D47').Value = '0.0₈116'
$ws.Range('E47').Value = '  +2.64%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.4561'
$ws.Range('E48').Value = '  -0.62%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '55.54'
$ws.Range('E49').Value = '  -4.26%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.008'
$ws.Range('E50').Value = '  -0.05%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.885'
$ws.Range('E51').Value = '  -3.66%  '
